# config and Test Case updated
#
# The "Run Mode" column (C) on the "Automation Tests" sheet stores a Yes/No
# flag per test case. A block of test cases that used to be enabled ("Yes")
# is being disabled, while the block of test cases that follows (previously
# disabled) is being enabled instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Automation Tests")
$ws.Activate()

# Rows that flip from "Yes" -> "No"
$rowsToNo = @(4, 32, 33, 34, 35)
foreach ($r in $rowsToNo) {
    $ws.Cells.Item($r, 3).Value = "No"
}

# Rows that flip from "No" -> "Yes"
$rowsToYes = 38..84
foreach ($r in $rowsToYes) {
    $ws.Cells.Item($r, 3).Value = "Yes"
}

# Update the saved view/selection to match the newly relevant block of rows.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C38:C84").Select()
